$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 8, pushing the existing "extr1".."extr8" rows
# down by two (they become rows 10..17), to make room for two new
# contingency lines ("line7" and "line8").
$ws.Rows.Item(8).Resize(2).Insert()

# New "line7" row (row 8)
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# New "line8" row (row 9)
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Match the column-A number formatting/style used by the other data rows
# (the Insert() above only shifted formatting from row 8 downward, leaving
# the two brand new rows without it).
$ws.Range("A10").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Update the "in_service" flags for the shifted extr rows.
$ws.Cells.Item(11, 5).Value = $true    # extr2
$ws.Cells.Item(12, 5).Value = $true    # extr3
$ws.Cells.Item(13, 5).Value = $true    # extr4
$ws.Cells.Item(14, 5).Value = $false   # extr5
$ws.Cells.Item(15, 5).Value = $true    # extr6
